$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.868.08'
$ws.Range("E2").Value = '  -1.05%  '

$ws.Range("D3").Value = '3.425.72'
$ws.Range("E3").Value = '  -2.41%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.11%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '3.427.26'
$ws.Range("E8").Value = '  -2.20%  '

$ws.Range("E9").Value = '  -5.62%  '

$ws.Range("E10").Value = '  -9.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.01'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.82%  '

$ws.Range("E12").Value = '  -6.71%  '

$ws.Range("D13").Value = '4.009.28'
$ws.Range("E13").Value = '  -2.48%  '

$ws.Range("E14").Value = '  -8.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.61%  '

$ws.Range("D16").Value = '3.411.31'
$ws.Range("E16").Value = '  -2.80%  '

$ws.Range("E17").Value = '  -1.85%  '

$ws.Range("D18").Value = '64.755.17'
$ws.Range("E18").Value = '  -1.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -11.94%  '

$ws.Range("E20").Value = '  -5.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '383.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.24%  '

$ws.Range("D26").Value = '3.563.50'
$ws.Range("E26").Value = '  -2.51%  '

$ws.Range("E27").Value = '  -6.86%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.59%  '

$ws.Range("E30").Value = '  -8.88%  '

$ws.Range("E31").Value = '  -10.81%  '

$ws.Range("D32").Value = '3.437.12'
$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.11%  '

$ws.Range("E35").Value = '  -7.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.87'
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = '  -11.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.93%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.62'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.61%  '

$ws.Range("E41").Value = '  -6.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.810'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -13.44%  '

$ws.Range("E46").Value = '  -7.78%  '

$ws.Range("E47").Value = '  +4.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.45'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.02'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -13.86%  '

$ws.Range("D51").Value = '2.165.72'
$ws.Range("E51").Value = '  -6.83%  '
